$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 38 with the new question entry
$ws.Range("A38").Value = 450
$ws.Range("B38").Value = "NA"
$ws.Range("C38").Value = "Delete Node in a BST"

# Match formatting/styles of the row above (row 37) for consistency
$ws.Range("A37").Copy()
$ws.Range("A38").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B37").Copy()
$ws.Range("B38").PasteSpecial(-4122)
$ws.Range("C37").Copy()
$ws.Range("C38").PasteSpecial(-4122)

$excel.CutCopyMode = 0
